$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 scores
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3

# Update row 16 scores
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = 3

# Update row 17 scores
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 3

# Update the view: frozen pane top-left cell and selection
$ws.Activate()
$ws.Range("L16").Select()
$excel.ActiveWindow.ScrollRow = 12
